$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column Q (2020) into the new column R (2021)
# for the header row and every data row, then fill in the new values.
$ws.Range("Q3:Q33").Copy()
$ws.Range("R3:R33").PasteSpecial(-4122)

$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 12.6
$ws.Range("R5").Value = 17.9
$ws.Range("R6").Value = 7.3
$ws.Range("R7").Value = 12.6
$ws.Range("R8").Value = 19.5
$ws.Range("R9").Value = 5.5
$ws.Range("R10").Value = 10.3
$ws.Range("R11").Value = 12.3
$ws.Range("R12").Value = 8.2
$ws.Range("R13").Value = 24.8
$ws.Range("R14").Value = 33.1
$ws.Range("R15").Value = 16.6
$ws.Range("R16").Value = 23.9
$ws.Range("R17").Value = 29.5
$ws.Range("R18").Value = 18.1
$ws.Range("R19").Value = 9.6
$ws.Range("R20").Value = 14.8
$ws.Range("R21").Value = 4.3
$ws.Range("R22").Value = 12.1
$ws.Range("R23").Value = 18.2
$ws.Range("R24").Value = 5.9
$ws.Range("R25").Value = 17.3
$ws.Range("R26").Value = 27.6
$ws.Range("R27").Value = 7.4
$ws.Range("R28").Value = 7.8
$ws.Range("R29").Value = 10.4
$ws.Range("R30").Value = 5.6
$ws.Range("R31").Value = 6.7
$ws.Range("R32").Value = 10.7
$ws.Range("R33").Value = 3

# Leave the selection where the author left it after editing the sheet.
$ws.Range("S4").Select()
